$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44160
$ws.Range("N2").Value = 19000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 19500
$ws.Range("S2").Value = 1500
$ws.Range("D4").Value = 45126
$ws.Range("M4").Value = 160
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 14375
$ws.Range("Q4").Value = "$/bandeja 10 kilos"
$ws.Range("S4").Value = 1438
$ws.Range("T4").Value = 10
$ws.Range("D5").Value = 45126
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 180
$ws.Range("N5").Value = 13000
$ws.Range("O5").Value = 13000
$ws.Range("P5").Value = 13000
$ws.Range("S5").Value = 1300
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 160
$ws.Range("N6").Value = 14000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 14375
$ws.Range("S6").Value = 1438
$ws.Range("D7").Value = 45125
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 180
$ws.Range("N7").Value = 13000
$ws.Range("O7").Value = 13000
$ws.Range("P7").Value = 13000
$ws.Range("S7").Value = 1300
$ws.Range("D8").Value = 44496
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 23000
$ws.Range("O8").Value = 24000
$ws.Range("P8").Value = 23500
$ws.Range("Q8").Value = "$/caja 12 kilos"
$ws.Range("S8").Value = 1958
$ws.Range("T8").Value = 12
$ws.Range("D9").Value = 44468
$ws.Range("M9").Value = 200
$ws.Range("Q9").Value = "$/bandeja 10 kilos"
$ws.Range("S9").Value = 2950
$ws.Range("T9").Value = 10
$ws.Range("D10").Value = 44524
$ws.Range("R10").Value = "Región de Coquimbo"
$ws.Range("D11").Value = 44860
$ws.Range("N11").Value = 23000
$ws.Range("O11").Value = 24000
$ws.Range("P11").Value = 23500
$ws.Range("R11").Value = "Provincia de Limarí"
$ws.Range("S11").Value = 1958
$ws.Range("D12").Value = 44839
$ws.Range("L12").Value = "Segunda"
$ws.Range("M12").Value = 160
$ws.Range("N12").Value = 26000
$ws.Range("O12").Value = 27000
$ws.Range("P12").Value = 26500
$ws.Range("Q12").Value = "$/caja 12 kilos"
$ws.Range("S12").Value = 2208
$ws.Range("D13").Value = 44489
$ws.Range("M13").Value = 200
$ws.Range("N13").Value = 24000
$ws.Range("O13").Value = 25000
$ws.Range("P13").Value = 24500
$ws.Range("Q13").Value = "$/caja 12 kilos"
$ws.Range("S13").Value = 2042
$ws.Range("T13").Value = 12
$ws.Range("D14").Value = 44441
$ws.Range("L14").Value = "Primera"
$ws.Range("N14").Value = 29000
$ws.Range("O14").Value = 30000
$ws.Range("P14").Value = 29500
$ws.Range("S14").Value = 2458
$ws.Range("D15").Value = 44776
$ws.Range("L15").Value = "Segunda"
$ws.Range("M15").Value = 160
$ws.Range("N15").Value = 29000
$ws.Range("O15").Value = 30000
$ws.Range("P15").Value = 29500
$ws.Range("Q15").Value = "$/caja 10 kilos"
$ws.Range("S15").Value = 2950
$ws.Range("T15").Value = 10
$ws.Range("D16").Value = 44811
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = 29000
$ws.Range("O16").Value = 30000
$ws.Range("P16").Value = 29500
$ws.Range("S16").Value = 2458
$ws.Range("D17").Value = 44475
$ws.Range("L17").Value = "Especial"
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = 32000
$ws.Range("O17").Value = 33000
$ws.Range("P17").Value = 32500
$ws.Range("Q17").Value = "$/caja 12 kilos"
$ws.Range("S17").Value = 2708
$ws.Range("T17").Value = 12
$ws.Range("D18").Value = 44881
$ws.Range("Q18").Value = "$/caja 12 kilos"
$ws.Range("S18").Value = 1875
$ws.Range("T18").Value = 12
$ws.Range("D19").Value = 44545
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 200
$ws.Range("N19").Value = 23000
$ws.Range("O19").Value = 24000
$ws.Range("P19").Value = 23500
$ws.Range("Q19").Value = "$/bandeja 12 kilos"
$ws.Range("S19").Value = 1958
$ws.Range("D20").Value = 45205
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = 22000
$ws.Range("O20").Value = 23000
$ws.Range("P20").Value = 22500
$ws.Range("Q20").Value = "$/bandeja 10 kilos"
$ws.Range("S20").Value = 2250
$ws.Range("T20").Value = 10
$ws.Range("D21").Value = 44874
$ws.Range("M21").Value = 250
$ws.Range("D22").Value = 45147
$ws.Range("M22").Value = 270
$ws.Range("N22").Value = 17000
$ws.Range("O22").Value = 18000
$ws.Range("P22").Value = 17500
$ws.Range("Q22").Value = "$/caja 10 kilos"
$ws.Range("S22").Value = 1750
$ws.Range("T22").Value = 10
$ws.Range("D23").Value = 44846
$ws.Range("M23").Value = 160
$ws.Range("N23").Value = 24000
$ws.Range("O23").Value = 25000
$ws.Range("P23").Value = 24500
$ws.Range("Q23").Value = "$/caja 12 kilos"
$ws.Range("S23").Value = 2042
$ws.Range("T23").Value = 12
$ws.Range("D24").Value = 44846
$ws.Range("L24").Value = "Segunda"
$ws.Range("M24").Value = 100
$ws.Range("N24").Value = 22000
$ws.Range("O24").Value = 23000
$ws.Range("P24").Value = 22500
$ws.Range("S24").Value = 1875
$ws.Range("D25").Value = 45133
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 150
$ws.Range("N25").Value = 15000
$ws.Range("O25").Value = 16000
$ws.Range("P25").Value = 15667
$ws.Range("Q25").Value = "$/bandeja 10 kilos"
$ws.Range("S25").Value = 1567
$ws.Range("T25").Value = 10
$ws.Range("D26").Value = 44167
$ws.Range("L26").Value = "Segunda"
$ws.Range("N26").Value = 18000
$ws.Range("O26").Value = 19000
$ws.Range("P26").Value = 18500
$ws.Range("Q26").Value = "$/caja 13 kilos"
$ws.Range("S26").Value = 1423
$ws.Range("T26").Value = 13
$ws.Range("D27").Value = 44783
$ws.Range("L27").Value = "Tercera"
$ws.Range("M27").Value = 100
$ws.Range("N27").Value = 27000
$ws.Range("O27").Value = 28000
$ws.Range("P27").Value = 27500
$ws.Range("S27").Value = 2292
